# Apply the "Updated symbol list" crypto price/coin-listing refresh.
# Values in column D look numeric but the source data stores them as plain text
# (inline strings); prefixing with a leading apostrophe forces Excel to keep them
# as text instead of re-typing the cell as a Number, matching the original file.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.21"
$ws.Range("D3").Value = "'24.03"
$ws.Range("D4").Value = "'5.362"
$ws.Range("D5").Value = "'0.05823"
$ws.Range("D6").Value = "'3.372"
$ws.Range("D7").Value = "'6.476"
$ws.Range("D8").Value = "'0.8100"
$ws.Range("D9").Value = "'0.9213"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1403"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07407"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03197"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03027"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09386"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.852"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001568"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04723"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005967"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006055"
$ws.Range("D21").Value = "'0.004690"
$ws.Range("D22").Value = "'0.00008797"
$ws.Range("E22").Value = "21NitroExNTXBestin24h"
$ws.Range("D23").Value = "'3.599"
$ws.Range("D25").Value = "'0.3179"
$ws.Range("D26").Value = "'0.1319"
$ws.Range("D28").Value = "'0.0002349"
$ws.Range("D40").Value = "'0.03853"
$ws.Range("D41").Value = "'0.006337"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1068"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002749"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.009067"
$ws.Range("D45").Value = "'0.00005262"
$ws.Range("D46").Value = "'0.00000000750"
$ws.Range("D47").Value = "'0.7097"
$ws.Range("D48").Value = "'0.001841"
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("D50").Value = "'0.0001999"

# The leading apostrophe above flips those cells to a "quote prefix" style so
# Excel remembers to redisplay them verbatim; reset the style back to Normal so
# no stray formatting is introduced (the source cells carry no explicit style).
$priceCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D21","D22","D23","D25","D26","D28","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50")
foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

